$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates ---
# E11: total "VALOR MORA" amount
$ws.Range("E11").Value = 549437

# C13: "Cant. Trabajadores" (worker count) 5 -> 4
$ws.Range("C13").Value = 4

# --- Remove the TOMAS IGNACIO ZAMBRANO BOLIVAR block (old rows 29:34) ---
# This shifts the VICTOR block + footer rows up by six, matching the
# target dimension (B2:J40) and footer rows moving from 45/46 to 39/40.
$ws.Rows("29:34").Delete()

# --- Rewrite the detail table (rows 16-34) with the updated data ---
# Clear old contents first so the shared-string pool is rebuilt in the
# same order the new values are entered below.
$ws.Range("C16:G34").ClearContents()

# Row 16: LILA MARINA ZAMBRANO BOLIVAR, periodo 1609
$ws.Range("C16").Value = "45534413"
$ws.Range("D16").Value = "LILA MARINA ZAMBRANO BOLIVAR"
$ws.Range("E16").Value = "1609"
$ws.Range("F16").Value = 27578
$ws.Range("G16").Value = 737717

# Rows 17-34: MELISSA / GUILLERMO / VICTOR, periodos 1901-1906
$workers = @(
  @{ Id = "1047454838"; Name = "MELISSA ANDREA COGOLLO MOLINA" },
  @{ Id = "1047365598"; Name = "GUILLERMO JOSE ZAMBRANO BOLIVAR" },
  @{ Id = "1047432113"; Name = "VICTOR ANDRES COGOLLO MOLINA" }
)
$periods = @("1901", "1902", "1903", "1904", "1905", "1906")

$row = 17
foreach ($periodo in $periods) {
  foreach ($worker in $workers) {
    $ws.Range("C$row").Value = $worker.Id
    $ws.Range("D$row").Value = $worker.Name
    $ws.Range("E$row").Value = $periodo
    if ($periodo -eq "1906") {
      $ws.Range("F$row").Value = 17708
    } else {
      $ws.Range("F$row").Value = 31249
    }
    $ws.Range("G$row").Value = 781242
    $row = $row + 1
  }
}
